# Insert a new weekly record at row 35 (shifting the existing rows 35-57
# down to 36-58) and populate it with the new data point.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Inserting the whole row pushes rows 35..57 down to 36..58 and copies the
# formatting from the row above (so D35 already inherits the date style).
$ws.Rows(35).Insert()

$ws.Range("A35").Value = 10
$ws.Range("B35").Value = "Vega Modelo de Temuco"
$ws.Range("C35").Value = "La Araucanía"
$ws.Range("D35").Value = 44777
$ws.Range("E35").Value = 9
$ws.Range("F35").Value = 100112010
$ws.Range("G35").Value = "Achicoria"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 65
$ws.Range("K35").Value = 10000
$ws.Range("L35").Value = 10000
$ws.Range("M35").Value = 10000
$ws.Range("N35").Value = "$/caja 18 unidades"
$ws.Range("O35").Value = "Región Metropolitana"
$ws.Range("P35").Value = 556
$ws.Range("Q35").Value = 18
$ws.Range("R35").Value = "Hortaliza"
